$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Name = "ROW50-FE-LIFTER";  Row = 57; A = 45752.71053879629; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x5e"; E = "0xe";  F = 400; G = [double]"5.68631262647114e+23"; GText = $null; H = 350; I = 14 },
    @{ Name = "ROW50-MID-LIFTER"; Row = 59; A = 45752.67828703704; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x62"; E = "0x19"; F = 400; G = 0; GText = "568631262647113771663628"; H = 354; I = 25 },
    @{ Name = "ROW11-FE-LIFTER";  Row = 57; A = 45752.74055164352; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x01,0x5e"; E = "0x14"; F = 400; G = [double]"5.68631262647114e+23"; GText = $null; H = 350; I = 20 },
    @{ Name = "ROW11-MID-LIFTER"; Row = 57; A = 45752.87703410879; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x01,0x66"; E = "0x19"; F = 400; G = [double]"5.68631262647114e+23"; GText = $null; H = 358; I = 25 }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)
    $r = $s.Row

    # Column A: date/time value, same display format as the row above it
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
    $ws.Cells.Item($r, 1).Value = $s.A

    # Columns B-E: text values (hex byte strings)
    $ws.Cells.Item($r, 2).Value = $s.B
    $ws.Cells.Item($r, 3).Value = $s.C
    $ws.Cells.Item($r, 4).Value = $s.D
    $ws.Cells.Item($r, 5).Value = $s.E

    # Column F: plain numeric value
    $ws.Cells.Item($r, 6).Value = $s.F

    # Column G: numeric value, except on ROW50-MID-LIFTER where it is stored as text
    if ($s.GText -ne $null) {
        $ws.Cells.Item($r, 7).NumberFormat = "@"
        $ws.Cells.Item($r, 7).Value = $s.GText
        $ws.Cells.Item($r, 7).ClearFormats()
    } else {
        $ws.Cells.Item($r, 7).Value = $s.G
    }

    # Columns H-I: plain numeric values
    $ws.Cells.Item($r, 8).Value = $s.H
    $ws.Cells.Item($r, 9).Value = $s.I
}
